$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2,1).Value = "FAPs"
$ws.Cells.Item(2,2).Value = "Lgi1"
$ws.Cells.Item(2,3).Value = "Adam11"
$ws.Cells.Item(2,4).Value = "ECs"
$ws.Cells.Item(2,5).Value = 2
$ws.Cells.Item(2,6).Value = 0.6666666666666666
$ws.Cells.Item(2,7).Value = 0.014884
$ws.Cells.Item(2,8).Value = 0.044652
$ws.Cells.Item(2,9).Value = 0.627584365205414
$ws.Cells.Item(2,10).Value = 0.627584365205414
$ws.Cells.Item(2,11).Value = 1
$ws.Cells.Item(2,12).Value = 0.3333333333333333
$ws.Cells.Item(2,13).Value = 0.075142
$ws.Cells.Item(2,14).Value = 0.225426
$ws.Cells.Item(2,15).Value = 0.05097354113547087
$ws.Cells.Item(2,16).Value = 0.05097354113547086
$ws.Cells.Item(2,17).Value = 0.001118413528
$ws.Cells.Item(2,18).Value = 0.010065721752
$ws.Cells.Item(2,19).Value = 0.03199019745577655
$ws.Cells.Item(2,20).Value = 0.03199019745577654

# Row 3
$ws.Cells.Item(3,1).Value = "FAPs"
$ws.Cells.Item(3,2).Value = "Lgi1"
$ws.Cells.Item(3,3).Value = "Adam11"
$ws.Cells.Item(3,4).Value = "FAPs"
$ws.Cells.Item(3,5).Value = 2
$ws.Cells.Item(3,6).Value = 0.6666666666666666
$ws.Cells.Item(3,7).Value = 0.014884
$ws.Cells.Item(3,8).Value = 0.044652
$ws.Cells.Item(3,9).Value = 0.627584365205414
$ws.Cells.Item(3,10).Value = 0.627584365205414
$ws.Cells.Item(3,11).Value = 3
$ws.Cells.Item(3,12).Value = 1
$ws.Cells.Item(3,13).Value = 0.5007510000000001
$ws.Cells.Item(3,14).Value = 1.502253
$ws.Cells.Item(3,15).Value = 0.3396908745725183
$ws.Cells.Item(3,16).Value = 0.3396908745725183
$ws.Cells.Item(3,17).Value = 0.007453177884000001
$ws.Cells.Item(3,18).Value = 0.067078600956
$ws.Cells.Item(3,19).Value = 0.2131846818846659
$ws.Cells.Item(3,20).Value = 0.2131846818846659

# Row 4
$ws.Cells.Item(4,1).Value = "FAPs"
$ws.Cells.Item(4,2).Value = "Lgi1"
$ws.Cells.Item(4,3).Value = "Adam11"
$ws.Cells.Item(4,4).Value = "Inflammatory-Mac"
$ws.Cells.Item(4,5).Value = 2
$ws.Cells.Item(4,6).Value = 0.6666666666666666
$ws.Cells.Item(4,7).Value = 0.014884
$ws.Cells.Item(4,8).Value = 0.044652
$ws.Cells.Item(4,9).Value = 0.627584365205414
$ws.Cells.Item(4,10).Value = 0.627584365205414
$ws.Cells.Item(4,11).Value = 2
$ws.Cells.Item(4,12).Value = 0.6666666666666666
$ws.Cells.Item(4,13).Value = 0.4588233333333334
$ws.Cells.Item(4,14).Value = 1.37647
$ws.Cells.Item(4,15).Value = 0.3112487031963553
$ws.Cells.Item(4,16).Value = 0.3112487031963553
$ws.Cells.Item(4,17).Value = 0.006829126493333334
$ws.Cells.Item(4,18).Value = 0.06146213844
$ws.Cells.Item(4,19).Value = 0.1953348198164929
$ws.Cells.Item(4,20).Value = 0.1953348198164929

# Row 5
$ws.Cells.Item(5,1).Value = "FAPs"
$ws.Cells.Item(5,2).Value = "Lgi1"
$ws.Cells.Item(5,3).Value = "Adam11"
$ws.Cells.Item(5,4).Value = "MuSCs"
$ws.Cells.Item(5,5).Value = 2
$ws.Cells.Item(5,6).Value = 0.6666666666666666
$ws.Cells.Item(5,7).Value = 0.014884
$ws.Cells.Item(5,8).Value = 0.044652
$ws.Cells.Item(5,9).Value = 0.627584365205414
$ws.Cells.Item(5,10).Value = 0.627584365205414
$ws.Cells.Item(5,11).Value = 3
$ws.Cells.Item(5,12).Value = 1
$ws.Cells.Item(5,13).Value = 0.4048450000000001
$ws.Cells.Item(5,14).Value = 1.214535
$ws.Cells.Item(5,15).Value = 0.2746318072581207
$ws.Cells.Item(5,16).Value = 0.2746318072581207
$ws.Cells.Item(5,17).Value = 0.006025712980000001
$ws.Cells.Item(5,18).Value = 0.05423141682
$ws.Cells.Item(5,19).Value = 0.1723546284233033
$ws.Cells.Item(5,20).Value = 0.1723546284233033

# Row 6
$ws.Cells.Item(6,1).Value = "FAPs"
$ws.Cells.Item(6,2).Value = "Lgi1"
$ws.Cells.Item(6,3).Value = "Adam11"
$ws.Cells.Item(6,4).Value = "Resolving-Mac"
$ws.Cells.Item(6,5).Value = 2
$ws.Cells.Item(6,6).Value = 0.6666666666666666
$ws.Cells.Item(6,7).Value = 0.014884
$ws.Cells.Item(6,8).Value = 0.044652
$ws.Cells.Item(6,9).Value = 0.627584365205414
$ws.Cells.Item(6,10).Value = 0.627584365205414
$ws.Cells.Item(6,11).Value = 1
$ws.Cells.Item(6,12).Value = 0.3333333333333333
$ws.Cells.Item(6,13).Value = 0.034576
$ws.Cells.Item(6,14).Value = 0.103728
$ws.Cells.Item(6,15).Value = 0.02345507383753481
$ws.Cells.Item(6,16).Value = 0.02345507383753481
$ws.Cells.Item(6,17).Value = 0.000514629184
$ws.Cells.Item(6,18).Value = 0.004631662656
$ws.Cells.Item(6,19).Value = 0.0147200376251754
$ws.Cells.Item(6,20).Value = 0.0147200376251754

# Row 7
$ws.Cells.Item(7,1).Value = "MuSCs"
$ws.Cells.Item(7,2).Value = "Lgi1"
$ws.Cells.Item(7,3).Value = "Adam11"
$ws.Cells.Item(7,4).Value = "ECs"
$ws.Cells.Item(7,5).Value = 2
$ws.Cells.Item(7,6).Value = 0.6666666666666666
$ws.Cells.Item(7,7).Value = 0.008832333333333333
$ws.Cells.Item(7,8).Value = 0.026497
$ws.Cells.Item(7,9).Value = 0.372415634794586
$ws.Cells.Item(7,10).Value = 0.3724156347945861
$ws.Cells.Item(7,11).Value = 1
$ws.Cells.Item(7,12).Value = 0.3333333333333333
$ws.Cells.Item(7,13).Value = 0.075142
$ws.Cells.Item(7,14).Value = 0.225426
$ws.Cells.Item(7,15).Value = 0.05097354113547087
$ws.Cells.Item(7,16).Value = 0.05097354113547086
$ws.Cells.Item(7,17).Value = 0.0006636791913333333
$ws.Cells.Item(7,18).Value = 0.005973112722
$ws.Cells.Item(7,19).Value = 0.01898334367969432
$ws.Cells.Item(7,20).Value = 0.01898334367969432

# Row 8
$ws.Cells.Item(8,1).Value = "MuSCs"
$ws.Cells.Item(8,2).Value = "Lgi1"
$ws.Cells.Item(8,3).Value = "Adam11"
$ws.Cells.Item(8,4).Value = "FAPs"
$ws.Cells.Item(8,5).Value = 2
$ws.Cells.Item(8,6).Value = 0.6666666666666666
$ws.Cells.Item(8,7).Value = 0.008832333333333333
$ws.Cells.Item(8,8).Value = 0.026497
$ws.Cells.Item(8,9).Value = 0.372415634794586
$ws.Cells.Item(8,10).Value = 0.3724156347945861
$ws.Cells.Item(8,11).Value = 3
$ws.Cells.Item(8,12).Value = 1
$ws.Cells.Item(8,13).Value = 0.5007510000000001
$ws.Cells.Item(8,14).Value = 1.502253
$ws.Cells.Item(8,15).Value = 0.3396908745725183
$ws.Cells.Item(8,16).Value = 0.3396908745725183
$ws.Cells.Item(8,17).Value = 0.004422799749
$ws.Cells.Item(8,18).Value = 0.039805197741
$ws.Cells.Item(8,19).Value = 0.1265061926878525
$ws.Cells.Item(8,20).Value = 0.1265061926878525

# Row 9
$ws.Cells.Item(9,1).Value = "MuSCs"
$ws.Cells.Item(9,2).Value = "Lgi1"
$ws.Cells.Item(9,3).Value = "Adam11"
$ws.Cells.Item(9,4).Value = "Inflammatory-Mac"
$ws.Cells.Item(9,5).Value = 2
$ws.Cells.Item(9,6).Value = 0.6666666666666666
$ws.Cells.Item(9,7).Value = 0.008832333333333333
$ws.Cells.Item(9,8).Value = 0.026497
$ws.Cells.Item(9,9).Value = 0.372415634794586
$ws.Cells.Item(9,10).Value = 0.3724156347945861
$ws.Cells.Item(9,11).Value = 2
$ws.Cells.Item(9,12).Value = 0.6666666666666666
$ws.Cells.Item(9,13).Value = 0.4588233333333334
$ws.Cells.Item(9,14).Value = 1.37647
$ws.Cells.Item(9,15).Value = 0.3112487031963553
$ws.Cells.Item(9,16).Value = 0.3112487031963553
$ws.Cells.Item(9,17).Value = 0.004052480621111111
$ws.Cells.Item(9,18).Value = 0.03647232559
$ws.Cells.Item(9,19).Value = 0.1159138833798623
$ws.Cells.Item(9,20).Value = 0.1159138833798624

# Row 10
$ws.Cells.Item(10,1).Value = "MuSCs"
$ws.Cells.Item(10,2).Value = "Lgi1"
$ws.Cells.Item(10,3).Value = "Adam11"
$ws.Cells.Item(10,4).Value = "MuSCs"
$ws.Cells.Item(10,5).Value = 2
$ws.Cells.Item(10,6).Value = 0.6666666666666666
$ws.Cells.Item(10,7).Value = 0.008832333333333333
$ws.Cells.Item(10,8).Value = 0.026497
$ws.Cells.Item(10,9).Value = 0.372415634794586
$ws.Cells.Item(10,10).Value = 0.3724156347945861
$ws.Cells.Item(10,11).Value = 3
$ws.Cells.Item(10,12).Value = 1
$ws.Cells.Item(10,13).Value = 0.4048450000000001
$ws.Cells.Item(10,14).Value = 1.214535
$ws.Cells.Item(10,15).Value = 0.2746318072581207
$ws.Cells.Item(10,16).Value = 0.2746318072581207
$ws.Cells.Item(10,17).Value = 0.003575725988333334
$ws.Cells.Item(10,18).Value = 0.032181533895
$ws.Cells.Item(10,19).Value = 0.1022771788348174
$ws.Cells.Item(10,20).Value = 0.1022771788348174

# Row 11
$ws.Cells.Item(11,1).Value = "MuSCs"
$ws.Cells.Item(11,2).Value = "Lgi1"
$ws.Cells.Item(11,3).Value = "Adam11"
$ws.Cells.Item(11,4).Value = "Resolving-Mac"
$ws.Cells.Item(11,5).Value = 2
$ws.Cells.Item(11,6).Value = 0.6666666666666666
$ws.Cells.Item(11,7).Value = 0.008832333333333333
$ws.Cells.Item(11,8).Value = 0.026497
$ws.Cells.Item(11,9).Value = 0.372415634794586
$ws.Cells.Item(11,10).Value = 0.3724156347945861
$ws.Cells.Item(11,11).Value = 1
$ws.Cells.Item(11,12).Value = 0.3333333333333333
$ws.Cells.Item(11,13).Value = 0.034576
$ws.Cells.Item(11,14).Value = 0.103728
$ws.Cells.Item(11,15).Value = 0.02345507383753481
$ws.Cells.Item(11,16).Value = 0.02345507383753481
$ws.Cells.Item(11,17).Value = 0.0003053867573333334
$ws.Cells.Item(11,18).Value = 0.002748480816
$ws.Cells.Item(11,19).Value = 0.008735036212359413
$ws.Cells.Item(11,20).Value = 0.008735036212359413
